$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$values = @{
    2 = 2270
    3 = 1714
    4 = 334
    5 = 1090
    6 = 824
    7 = 39
    8 = 5837
    9 = 89
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Range("F$row").Value = $values[$row]
    }
}
